# Apply the "Zhishang's inheritance" edit:
#  - The sub-header in I3 (under "Assigned to which Continuing Member?") is
#    rephrased from describing the *receiving* person to describing the
#    *passing* person.
#  - Every row's "Assigned to which Continuing Member?" cell (I4:I13) is
#    updated to the name of the person handing the project over: Zhishang.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = "Name of the person passing the inheritance."

$ws.Range("I4").Value = "Zhishang"
$ws.Range("I5").Value = "Zhishang"
$ws.Range("I6").Value = "Zhishang"
$ws.Range("I7").Value = "Zhishang"
$ws.Range("I8").Value = "Zhishang"
$ws.Range("I9").Value = "Zhishang"
$ws.Range("I10").Value = "Zhishang"
$ws.Range("I11").Value = "Zhishang"
$ws.Range("I12").Value = "Zhishang"
$ws.Range("I13").Value = "Zhishang"
